# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Wed Jan 10 17:10:19 UTC 2024 with GitHub Actions".
# Only the D (Price) and E (Volume(1h)) columns change; everything else is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.342.78'
$ws.Range("E2").Value = '  -3.25%  '
$ws.Range("D3").Value = '2.410.57'
$ws.Range("E3").Value = '  +6.63%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.73%  '
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.500'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0778'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").Value = '2.774.49'
$ws.Range("E14").Value = '  +6.45%  '
$ws.Range("D15").Value = '2.415.36'
$ws.Range("E15").Value = '  +6.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.834'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.05%  '
$ws.Range("D18").Value = '45.215.65'
$ws.Range("E18").Value = '  -3.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.80%  '
$ws.Range("D20").Value = '0.0₃0939'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.49%  '
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.48%  '
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '148.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +17.93%  '
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0297'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = '1.993.44'
$ws.Range("E42").Value = '  +11.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.79'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +25.12%  '
$ws.Range("E48").Value = '  +9.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.53'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.99%  '
$ws.Range("D50").Value = '2.644.30'
$ws.Range("E50").Value = '  +6.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.182'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.15%  '
